# chore: update Sheets via scheduled runner
# Refreshes cached Universalis market-board averages (currentAveragePrice*)
# and the derived Leve-profit figures (LevePrice*/LeveProfit*) for the rows
# whose source prices moved since the last scheduled pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce
$ws.Cells.Item(11, 8).Value = 52925.332  # H11
$ws.Cells.Item(11, 9).Value = 52925.332  # I11
$ws.Cells.Item(11, 11).Value = 52925.332  # K11
$ws.Cells.Item(11, 13).Value = -52785.332  # M11

# Row 38: Just Give Him a Serum
$ws.Cells.Item(38, 8).Value = 1203.3334  # H38
$ws.Cells.Item(38, 9).Value = 103.75  # I38
$ws.Cells.Item(38, 10).Value = 10000  # J38
$ws.Cells.Item(38, 11).Value = 311.25  # K38
$ws.Cells.Item(38, 12).Value = 30000  # L38
$ws.Cells.Item(38, 13).Value = 60.75  # M38
$ws.Cells.Item(38, 14).Value = -30744  # N38

# Row 40: Stuck in the Moment
$ws.Cells.Item(40, 9).Value = 2477.6667  # I40
$ws.Cells.Item(40, 10).Value = 2781.7273  # J40
$ws.Cells.Item(40, 11).Value = 2477.6667  # K40
$ws.Cells.Item(40, 12).Value = 2781.7273  # L40
$ws.Cells.Item(40, 13).Value = -2302.6667  # M40
$ws.Cells.Item(40, 14).Value = -3131.7273  # N40

# Row 42: Eye of the Beholder
$ws.Cells.Item(42, 8).Value = 89  # H42
$ws.Cells.Item(42, 9).Value = 23.333334  # I42
$ws.Cells.Item(42, 10).Value = 117.14286  # J42
$ws.Cells.Item(42, 11).Value = 70.00000199999999  # K42
$ws.Cells.Item(42, 12).Value = 351.42858  # L42
$ws.Cells.Item(42, 13).Value = 159.999998  # M42
$ws.Cells.Item(42, 14).Value = -811.42858  # N42

# Row 64: Forged from the Void
$ws.Cells.Item(64, 8).Value = 4178.2856  # H64
$ws.Cells.Item(64, 9).Value = 3749.3333  # I64
$ws.Cells.Item(64, 10).Value = 4500  # J64
$ws.Cells.Item(64, 11).Value = 3749.3333  # K64
$ws.Cells.Item(64, 12).Value = 4500  # L64
$ws.Cells.Item(64, 13).Value = -3501.3333  # M64
$ws.Cells.Item(64, 14).Value = -4996  # N64

# Row 67: Dodging the Draft (L)
$ws.Cells.Item(67, 8).Value = 4178.2856  # H67
$ws.Cells.Item(67, 9).Value = 3749.3333  # I67
$ws.Cells.Item(67, 10).Value = 4500  # J67
$ws.Cells.Item(67, 11).Value = 3749.3333  # K67
$ws.Cells.Item(67, 12).Value = 4500  # L67
$ws.Cells.Item(67, 13).Value = -2891.3333  # M67
$ws.Cells.Item(67, 14).Value = -6216  # N67

# Row 82: Rolling on Initiative
$ws.Cells.Item(82, 8).Value = 268.2  # H82
$ws.Cells.Item(82, 9).Value = 268.2  # I82
$ws.Cells.Item(82, 11).Value = 804.5999999999999  # K82
$ws.Cells.Item(82, 13).Value = -398.5999999999999  # M82

# Row 85: Darkly Dreaming Dexterity (L)
$ws.Cells.Item(85, 8).Value = 268.2  # H85
$ws.Cells.Item(85, 9).Value = 268.2  # I85
$ws.Cells.Item(85, 11).Value = 804.5999999999999  # K85
$ws.Cells.Item(85, 13).Value = 599.4000000000001  # M85

# Row 100: Asking for a Friend
$ws.Cells.Item(100, 8).Value = 5233.3335  # H100
$ws.Cells.Item(100, 9).Value = 6600.5557  # I100
$ws.Cells.Item(100, 10).Value = 1131.6666  # J100
$ws.Cells.Item(100, 11).Value = 6600.5557  # K100
$ws.Cells.Item(100, 12).Value = 1131.6666  # L100
$ws.Cells.Item(100, 13).Value = -6059.5557  # M100
$ws.Cells.Item(100, 14).Value = -2213.6666  # N100

# Row 135: For Tired Minds
$ws.Cells.Item(135, 8).Value = 37038492  # H135
$ws.Cells.Item(135, 9).Value = 40001508  # I135
$ws.Cells.Item(135, 10).Value = 780  # J135
$ws.Cells.Item(135, 11).Value = 360013572  # K135
$ws.Cells.Item(135, 12).Value = 7020  # L135
$ws.Cells.Item(135, 13).Value = -360011037  # M135
$ws.Cells.Item(135, 14).Value = -12090  # N135

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 2939.1292  # H138
$ws.Cells.Item(138, 9).Value = 1166.8695  # I138
$ws.Cells.Item(138, 11).Value = 3500.6085  # K138
$ws.Cells.Item(138, 13).Value = 1639.3915  # M138

$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector
$ws.Cells.Item(82, 8).Value = 46359.875  # H82
$ws.Cells.Item(82, 9).Value = 16966.666  # I82
$ws.Cells.Item(82, 10).Value = 63995.8  # J82
$ws.Cells.Item(82, 11).Value = 16966.666  # K82
$ws.Cells.Item(82, 12).Value = 63995.8  # L82
$ws.Cells.Item(82, 13).Value = -16583.666  # M82
$ws.Cells.Item(82, 14).Value = -64761.8  # N82

# Row 85: The Clamor for Hammers (L)
$ws.Cells.Item(85, 8).Value = 46359.875  # H85
$ws.Cells.Item(85, 9).Value = 16966.666  # I85
$ws.Cells.Item(85, 10).Value = 63995.8  # J85
$ws.Cells.Item(85, 11).Value = 16966.666  # K85
$ws.Cells.Item(85, 12).Value = 63995.8  # L85
$ws.Cells.Item(85, 13).Value = -15640.666  # M85
$ws.Cells.Item(85, 14).Value = -66647.8  # N85

# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 3893.7334  # H86
$ws.Cells.Item(86, 9).Value = 5314.75  # I86
$ws.Cells.Item(86, 10).Value = 2269.7144  # J86
$ws.Cells.Item(86, 11).Value = 5314.75  # K86
$ws.Cells.Item(86, 12).Value = 2269.7144  # L86
$ws.Cells.Item(86, 13).Value = -4191.75  # M86
$ws.Cells.Item(86, 14).Value = -4515.7144  # N86

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 3893.7334  # H89
$ws.Cells.Item(89, 9).Value = 5314.75  # I89
$ws.Cells.Item(89, 10).Value = 2269.7144  # J89
$ws.Cells.Item(89, 11).Value = 26573.75  # K89
$ws.Cells.Item(89, 12).Value = 11348.572  # L89
$ws.Cells.Item(89, 13).Value = -20957.75  # M89
$ws.Cells.Item(89, 14).Value = -22580.572  # N89

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 2866.627  # H31
$ws.Cells.Item(31, 9).Value = 2243.375  # I31
$ws.Cells.Item(31, 10).Value = 3214.4883  # J31
$ws.Cells.Item(31, 11).Value = 2243.375  # K31
$ws.Cells.Item(31, 12).Value = 3214.4883  # L31
$ws.Cells.Item(31, 13).Value = -1948.375  # M31
$ws.Cells.Item(31, 14).Value = -3804.4883  # N31

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 2866.627  # H34
$ws.Cells.Item(34, 9).Value = 2243.375  # I34
$ws.Cells.Item(34, 10).Value = 3214.4883  # J34
$ws.Cells.Item(34, 11).Value = 2243.375  # K34
$ws.Cells.Item(34, 12).Value = 3214.4883  # L34
$ws.Cells.Item(34, 13).Value = -2041.375  # M34
$ws.Cells.Item(34, 14).Value = -3618.4883  # N34

# Row 86: Birch, Please
$ws.Cells.Item(86, 8).Value = 4512.4287  # H86
$ws.Cells.Item(86, 9).Value = 3659.6  # I86
$ws.Cells.Item(86, 10).Value = 6644.5  # J86
$ws.Cells.Item(86, 11).Value = 3659.6  # K86
$ws.Cells.Item(86, 12).Value = 6644.5  # L86
$ws.Cells.Item(86, 13).Value = -2536.6  # M86
$ws.Cells.Item(86, 14).Value = -8890.5  # N86

# Row 89: Built This City on Blocks and Soul (L)
$ws.Cells.Item(89, 8).Value = 4512.4287  # H89
$ws.Cells.Item(89, 9).Value = 3659.6  # I89
$ws.Cells.Item(89, 10).Value = 6644.5  # J89
$ws.Cells.Item(89, 11).Value = 18298  # K89
$ws.Cells.Item(89, 12).Value = 33222.5  # L89
$ws.Cells.Item(89, 13).Value = -12682  # M89
$ws.Cells.Item(89, 14).Value = -44454.5  # N89

# Row 107: Built to Last
$ws.Cells.Item(107, 8).Value = 3217.0625  # H107
$ws.Cells.Item(107, 9).Value = 2709.0908  # I107
$ws.Cells.Item(107, 10).Value = 4334.6  # J107
$ws.Cells.Item(107, 11).Value = 2709.0908  # K107
$ws.Cells.Item(107, 12).Value = 4334.6  # L107
$ws.Cells.Item(107, 13).Value = -789.0907999999999  # M107
$ws.Cells.Item(107, 14).Value = -8174.6  # N107

# Row 134: Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 2616.8  # H134
$ws.Cells.Item(134, 9).Value = 2100.2  # I134
$ws.Cells.Item(134, 11).Value = 6300.599999999999  # K134
$ws.Cells.Item(134, 13).Value = -3765.599999999999  # M134

$ws = $wb.Worksheets.Item("CUL")
# Row 17: Chew the Fat
$ws.Cells.Item(17, 8).Value = 704.2727  # H17
$ws.Cells.Item(17, 9).Value = 456.16666  # I17
$ws.Cells.Item(17, 10).Value = 1002  # J17
$ws.Cells.Item(17, 11).Value = 1368.49998  # K17
$ws.Cells.Item(17, 12).Value = 3006  # L17
$ws.Cells.Item(17, 13).Value = -1199.49998  # M17
$ws.Cells.Item(17, 14).Value = -3344  # N17

# Row 22: A Total Nut Job
$ws.Cells.Item(22, 8).Value = 563  # H22
$ws.Cells.Item(22, 9).Value = 501  # I22
$ws.Cells.Item(22, 10).Value = 625  # J22
$ws.Cells.Item(22, 11).Value = 1503  # K22
$ws.Cells.Item(22, 12).Value = 1875  # L22
$ws.Cells.Item(22, 13).Value = -1334  # M22
$ws.Cells.Item(22, 14).Value = -2213  # N22

# Row 27: Brain Food
$ws.Cells.Item(27, 8).Value = 563  # H27
$ws.Cells.Item(27, 9).Value = 501  # I27
$ws.Cells.Item(27, 10).Value = 625  # J27
$ws.Cells.Item(27, 11).Value = 1503  # K27
$ws.Cells.Item(27, 12).Value = 1875  # L27
$ws.Cells.Item(27, 13).Value = -1401  # M27
$ws.Cells.Item(27, 14).Value = -2079  # N27

# Row 41: Gegeruju Gets Down
$ws.Cells.Item(41, 8).Value = 225  # H41
$ws.Cells.Item(41, 9).Value = 150  # I41
$ws.Cells.Item(41, 10).Value = 300  # J41
$ws.Cells.Item(41, 11).Value = 450  # K41
$ws.Cells.Item(41, 12).Value = 900  # L41
$ws.Cells.Item(41, 13).Value = -112  # M41
$ws.Cells.Item(41, 14).Value = -1576  # N41

# Row 44: No More Dumpster Diving
$ws.Cells.Item(44, 8).Value = 3049  # H44
$ws.Cells.Item(44, 9).Value = 599  # I44
$ws.Cells.Item(44, 11).Value = 1797  # K44
$ws.Cells.Item(44, 13).Value = -1399  # M44

# Row 55: Pagan Pastries
$ws.Cells.Item(55, 8).Value = 903.125  # H55
$ws.Cells.Item(55, 9).Value = 245.83333  # I55
$ws.Cells.Item(55, 10).Value = 2875  # J55
$ws.Cells.Item(55, 11).Value = 737.49999  # K55
$ws.Cells.Item(55, 12).Value = 8625  # L55
$ws.Cells.Item(55, 13).Value = -560.49999  # M55
$ws.Cells.Item(55, 14).Value = -8979  # N55

$ws = $wb.Worksheets.Item("GSM")
# Row 55: If You've Got It, Flaunt It
$ws.Cells.Item(55, 8).Value = 6003.4287  # H55
$ws.Cells.Item(55, 9).Value = 3670.6667  # I55
$ws.Cells.Item(55, 11).Value = 3670.6667  # K55
$ws.Cells.Item(55, 13).Value = -3343.6667  # M55

# Row 62: The Goggles, They Do Naught
$ws.Cells.Item(62, 14).ClearContents()  # N62 removed
$ws.Cells.Item(62, 8).Value = 47166.668  # H62
$ws.Cells.Item(62, 9).Value = 47166.668  # I62
$ws.Cells.Item(62, 10).Value = 0  # J62
$ws.Cells.Item(62, 11).Value = 47166.668  # K62
$ws.Cells.Item(62, 12).Value = 0  # L62
$ws.Cells.Item(62, 13).Value = -46480.668  # M62

# Row 65: Peril Never Wore Safety Goggles (L)
$ws.Cells.Item(65, 14).ClearContents()  # N65 removed
$ws.Cells.Item(65, 8).Value = 47166.668  # H65
$ws.Cells.Item(65, 9).Value = 47166.668  # I65
$ws.Cells.Item(65, 10).Value = 0  # J65
$ws.Cells.Item(65, 11).Value = 141500.004  # K65
$ws.Cells.Item(65, 12).Value = 0  # L65
$ws.Cells.Item(65, 13).Value = -138068.004  # M65

# Row 96: Bracelet for Impact
$ws.Cells.Item(96, 8).Value = 42263.332  # H96
$ws.Cells.Item(96, 10).Value = 42263.332  # J96
$ws.Cells.Item(96, 12).Value = 42263.332  # L96
$ws.Cells.Item(96, 14).Value = -47755.332  # N96

# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 2201  # H102
$ws.Cells.Item(102, 9).Value = 1823.6666  # I102
$ws.Cells.Item(102, 11).Value = 1823.6666  # K102
$ws.Cells.Item(102, 13).Value = -201.6666  # M102

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Cells.Item(7, 8).Value = 1954.4  # H7
$ws.Cells.Item(7, 9).Value = 1996.5  # I7
$ws.Cells.Item(7, 10).Value = 1891.25  # J7
$ws.Cells.Item(7, 11).Value = 1996.5  # K7
$ws.Cells.Item(7, 12).Value = 1891.25  # L7
$ws.Cells.Item(7, 13).Value = -1884.5  # M7
$ws.Cells.Item(7, 14).Value = -2115.25  # N7

# Row 46: Supply Side Logic
$ws.Cells.Item(46, 8).Value = 2508.7222  # H46
$ws.Cells.Item(46, 9).Value = 975.8333  # I46
$ws.Cells.Item(46, 10).Value = 3275.1667  # J46
$ws.Cells.Item(46, 11).Value = 975.8333  # K46
$ws.Cells.Item(46, 12).Value = 3275.1667  # L46
$ws.Cells.Item(46, 13).Value = -787.8333  # M46
$ws.Cells.Item(46, 14).Value = -3651.1667  # N46

# Row 55: It's Not a Job, It's a Calling
$ws.Cells.Item(55, 8).Value = 453  # H55
$ws.Cells.Item(55, 9).Value = 464.63635  # I55
$ws.Cells.Item(55, 10).Value = 410.33334  # J55
$ws.Cells.Item(55, 11).Value = 464.63635  # K55
$ws.Cells.Item(55, 12).Value = 410.33334  # L55
$ws.Cells.Item(55, 13).Value = -291.63635  # M55
$ws.Cells.Item(55, 14).Value = -756.33334  # N55

# Row 61: Spelling Me Softly
$ws.Cells.Item(61, 8).Value = 10944.556  # H61
$ws.Cells.Item(61, 9).Value = 9846.450000000001  # I61
$ws.Cells.Item(61, 10).Value = 14082  # J61
$ws.Cells.Item(61, 11).Value = 9846.450000000001  # K61
$ws.Cells.Item(61, 12).Value = 14082  # L61
$ws.Cells.Item(61, 13).Value = -9644.450000000001  # M61
$ws.Cells.Item(61, 14).Value = -14486  # N61

# Row 113: Peace in Rest
$ws.Cells.Item(113, 8).Value = 10944.556  # H113
$ws.Cells.Item(113, 9).Value = 9846.450000000001  # I113
$ws.Cells.Item(113, 10).Value = 14082  # J113
$ws.Cells.Item(113, 11).Value = 9846.450000000001  # K113
$ws.Cells.Item(113, 12).Value = 14082  # L113
$ws.Cells.Item(113, 13).Value = -7676.450000000001  # M113
$ws.Cells.Item(113, 14).Value = -18422  # N113

# Row 126: Battered Books
$ws.Cells.Item(126, 8).Value = 1954.4  # H126
$ws.Cells.Item(126, 9).Value = 1996.5  # I126
$ws.Cells.Item(126, 10).Value = 1891.25  # J126
$ws.Cells.Item(126, 11).Value = 5989.5  # K126
$ws.Cells.Item(126, 12).Value = 5673.75  # L126
$ws.Cells.Item(126, 13).Value = -3519.5  # M126
$ws.Cells.Item(126, 14).Value = -10613.75  # N126

# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 5768.3667  # H132
$ws.Cells.Item(132, 9).Value = 3447.9375  # I132
$ws.Cells.Item(132, 10).Value = 8420.286  # J132
$ws.Cells.Item(132, 11).Value = 10343.8125  # K132
$ws.Cells.Item(132, 12).Value = 25260.858  # L132
$ws.Cells.Item(132, 13).Value = -7813.8125  # M132
$ws.Cells.Item(132, 14).Value = -30320.858  # N132

# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 2309.6667  # H136
$ws.Cells.Item(136, 9).Value = 2135.0476  # I136
$ws.Cells.Item(136, 10).Value = 3532  # J136
$ws.Cells.Item(136, 11).Value = 6405.1428  # K136
$ws.Cells.Item(136, 12).Value = 10596  # L136
$ws.Cells.Item(136, 13).Value = -3855.1428  # M136
$ws.Cells.Item(136, 14).Value = -15696  # N136

$ws = $wb.Worksheets.Item("WVR")
# Row 32: Piling It On
$ws.Cells.Item(32, 8).Value = 24687.375  # H32
$ws.Cells.Item(32, 9).Value = 2499  # I32
$ws.Cells.Item(32, 10).Value = 27857.143  # J32
$ws.Cells.Item(32, 11).Value = 2499  # K32
$ws.Cells.Item(32, 12).Value = 27857.143  # L32
$ws.Cells.Item(32, 13).Value = -2182  # M32
$ws.Cells.Item(32, 14).Value = -28491.143  # N32

# Row 81: Where the Dragonflies, the Net Catches
$ws.Cells.Item(81, 8).Value = 5409425  # H81
$ws.Cells.Item(81, 9).Value = 1771.6296  # I81
$ws.Cells.Item(81, 10).Value = 20010090  # J81
$ws.Cells.Item(81, 11).Value = 3543.2592  # K81
$ws.Cells.Item(81, 12).Value = 40020180  # L81
$ws.Cells.Item(81, 13).Value = -2482.2592  # M81
$ws.Cells.Item(81, 14).Value = -40022302  # N81

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Cells.Item(84, 8).Value = 5409425  # H84
$ws.Cells.Item(84, 9).Value = 1771.6296  # I84
$ws.Cells.Item(84, 10).Value = 20010090  # J84
$ws.Cells.Item(84, 11).Value = 17716.296  # K84
$ws.Cells.Item(84, 12).Value = 200100900  # L84
$ws.Cells.Item(84, 13).Value = -12412.296  # M84
$ws.Cells.Item(84, 14).Value = -200111508  # N84

# Row 126: A Polished Purchase
$ws.Cells.Item(126, 8).Value = 2021.1  # H126
$ws.Cells.Item(126, 9).Value = 1759.5  # I126
$ws.Cells.Item(126, 11).Value = 5278.5  # K126
$ws.Cells.Item(126, 13).Value = -2808.5  # M126

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 2976.8113  # H132
$ws.Cells.Item(132, 9).Value = 2590.739  # I132
$ws.Cells.Item(132, 11).Value = 7772.217000000001  # K132
$ws.Cells.Item(132, 13).Value = -5242.217000000001  # M132
